# Add function to create excel with multiple sheets.
#
# This form used to carry two tabs: "Summary" (the header/cover sheet) and
# "Sheet1" (a per-group detail template used by the code that programmatically
# clones/fills it for each group). Going forward the single remaining sheet is
# renamed to "Description" and the old "Sheet1" template tab is removed from
# this workbook (its layout already lives on, and is driven from, the first
# sheet).

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Drop the old "Sheet1" detail-template tab.
$wb.Worksheets.Item("Sheet1").Delete()

# The remaining "Summary" tab becomes "Description".
$wb.Worksheets.Item("Summary").Name = "Description"

# Keep the Print_Titles defined name (rows 1:5 repeat on every printed page)
# pointing at the renamed sheet.
$wb.Names.Item("Description!Print_Titles").RefersTo = "=Description!`$1:`$5"
